$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: TRADING_ATTEMPT
$ws.Range("A2").Value = "2025-09-10T13:49:26.889804"
$ws.Range("B2").Value = "TRADING_ATTEMPT"
$ws.Range("C2").Value = "BTC"
$ws.Range("D2").Value = "UNKNOWN"
$ws.Range("E2").Value = 113902.0791768574
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = "ATTEMPT"
$ws.Range("L2").Value = "Attempting trade 1/1"

# Row 3: POSITION_OPENED
$ws.Range("A3").Value = "2025-09-10T13:49:28.786247"
$ws.Range("B3").Value = "POSITION_OPENED"
$ws.Range("C3").Value = "BTC"
$ws.Range("D3").Value = "UNKNOWN"
$ws.Range("E3").Value = 113902.0791768574
$ws.Range("F3").Value = 9600
$ws.Range("G3").Value = 40
$ws.Range("H3").Value = 1.714390973461102
$ws.Range("I3").Value = ""
$ws.Range("J3").Value = ""
$ws.Range("K3").Value = "SUCCESS"
$ws.Range("L3").Value = ""
